# Update DM integration fixture hierarchies
# Applies updated UUIDs and column-A widths across the four sheets of the
# typed-dimensions-2018-1 workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet: CodeSchemes -------------------------------------------------
$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")
$wsCodeSchemes.Columns(1).ColumnWidth = 27.885714285714286
$wsCodeSchemes.Range("A2").Value = "bdd1aac0-ad7a-4e5f-8442-c493f0167c24"

# --- Sheet: Codes ---------------------------------------------------------
$wsCodes = $wb.Worksheets.Item("Codes")
$wsCodes.Range("A2").Value = "833dda2a-145f-439f-9e26-7235914121fa"
$wsCodes.Range("A3").Value = "2ecb13a0-83b3-442e-86a3-01549ed7d78d"

# --- Sheet: Extensions ----------------------------------------------------
$wsExtensions = $wb.Worksheets.Item("Extensions")
$wsExtensions.Columns(1).ColumnWidth = 33.385714285714286
$wsExtensions.Range("A2").Value = "1fb961fc-4917-4978-936d-0a58dd65ee3d"

# --- Sheet: Members_dpmDimension ------------------------------------------
$wsMembers = $wb.Worksheets.Item("Members_dpmDimension")
$wsMembers.Columns(1).ColumnWidth = 30.08571428571429
$wsMembers.Range("A2").Value = "7143cf0c-e938-4c64-8cd5-a9f13314994b"
$wsMembers.Range("A3").Value = "ec5f13d0-51ff-43b5-95aa-d0fbfef07552"
